$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new "Codigo" column at the front (shifts Nombre..Stock right by one,
# B..J -> C..K). This matches the commit: "incluyendo el codigo producto".
# ---------------------------------------------------------------------------
$ws.Columns("A:A").Insert()

# Header (row 1): copy the header format (fill/border/font) onto the new A1
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "Codigo"

# Body cells: copy the plain body format onto the new column A cells
$ws.Range("B2").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)
$ws.Range("A2").Value = "-"
$ws.Range("A3").Value = "KORI1u66LIsBK"
$ws.Range("A4").Value = "KORIHMXVOy6hx"

# Product names for the newly-added rows (as entered by the importer)
$ws.Range("B3").Value = "120TH CENTURY BOYS 01"
$ws.Range("B4").Value = "2ATTACK ON TITAN VOL. 1"

# ---------------------------------------------------------------------------
# New "Stock" column (K) - every freshly-imported row starts at stock = 1
# ---------------------------------------------------------------------------
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "Stock"

$ws.Range("J2").Copy()
$ws.Range("K2:K4").PasteSpecial(-4122)
$ws.Range("K2").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("K4").Value = 1

# ---------------------------------------------------------------------------
# Column widths - widen the columns that needed more room for the new data
# (Codigo, Nombre, Descripcion, Precio, Serie); the rest keep their existing
# auto-fit widths (just shifted one column to the right).
# ---------------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 16
$ws.Columns("B:B").ColumnWidth = 25.333333333333336
$ws.Columns("C:C").ColumnWidth = 29.333333333333336
$ws.Columns("D:D").ColumnWidth = 12.166666666666666
$ws.Columns("F:F").ColumnWidth = 50.666666666666664

# ---------------------------------------------------------------------------
# Selection left on B2 after the import edits
# ---------------------------------------------------------------------------
$ws.Range("B2").Select()
